# "best results for 1 and 9"
# Updates the probability grid in A1:E10 on Sheet1 with the corrected
# (final) values from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = 0.0008486049627111905
$ws.Cells.Item(1, 2).Value2 = 0.9991513950372889
$ws.Cells.Item(1, 3).Value2 = 0.0008486049627111867
$ws.Cells.Item(1, 4).Value2 = 0.9991513950372889
$ws.Cells.Item(1, 5).Value2 = 0.9991513950372889
$ws.Cells.Item(2, 1).Value2 = 0.9991513950372889
$ws.Cells.Item(2, 2).Value2 = 0.0008486049627111873
$ws.Cells.Item(2, 3).Value2 = 0.9991513950372889
$ws.Cells.Item(2, 4).Value2 = 0.0008486049627111882
$ws.Cells.Item(2, 5).Value2 = 0.0008486049627111882
$ws.Cells.Item(3, 1).Value2 = 0.0007723126728216415
$ws.Cells.Item(3, 2).Value2 = 0.9999999823656832
$ws.Cells.Item(3, 3).Value2 = 0.07719490164828074
$ws.Cells.Item(3, 4).Value2 = 0.9999850752682813
$ws.Cells.Item(3, 5).Value2 = 0.014766012930259687
$ws.Cells.Item(4, 1).Value2 = 0.001967986922999256
$ws.Cells.Item(4, 2).Value2 = 0.0598665017412413
$ws.Cells.Item(4, 3).Value2 = 0.0032102021747555274
$ws.Cells.Item(4, 4).Value2 = 0.0017142080046726037
$ws.Cells.Item(4, 5).Value2 = 0.012204457330757168
$ws.Cells.Item(5, 1).Value2 = 0.00019088910243797146
$ws.Cells.Item(5, 2).Value2 = 0.00024889230184229396
$ws.Cells.Item(5, 3).Value2 = 0.0000003699529291226178
$ws.Cells.Item(5, 4).Value2 = 0.0000002738387308158181
$ws.Cells.Item(5, 5).Value2 = 0.0000000031800355378368657
$ws.Cells.Item(6, 1).Value2 = 0.999988690751936
$ws.Cells.Item(6, 2).Value2 = 0.9889318052133065
$ws.Cells.Item(6, 3).Value2 = 0.0000002624318536963093
$ws.Cells.Item(6, 4).Value2 = 0.0019229309587742243
$ws.Cells.Item(6, 5).Value2 = 0.9998750013764924
$ws.Cells.Item(7, 1).Value2 = 0.0012689150783167067
$ws.Cells.Item(7, 2).Value2 = 0.5938118169228063
$ws.Cells.Item(7, 3).Value2 = 0.9979657785344246
$ws.Cells.Item(7, 4).Value2 = 0.03053333798289417
$ws.Cells.Item(7, 5).Value2 = 0.9114628073374834
$ws.Cells.Item(8, 1).Value2 = 0.9999999995439988
$ws.Cells.Item(8, 2).Value2 = 0.00000004436475871691599
$ws.Cells.Item(8, 3).Value2 = 0.9656492121538018
$ws.Cells.Item(8, 4).Value2 = 0.08069122619316081
$ws.Cells.Item(8, 5).Value2 = 0.00002837932168242639
$ws.Cells.Item(9, 1).Value2 = 0.9171000930773086
$ws.Cells.Item(9, 2).Value2 = 0.9990842735241952
$ws.Cells.Item(9, 3).Value2 = 0.278858128312238
$ws.Cells.Item(9, 4).Value2 = 0.9999999899681725
$ws.Cells.Item(9, 5).Value2 = 0.9999937147734456
$ws.Cells.Item(10, 1).Value2 = 0.014441997333328998
$ws.Cells.Item(10, 2).Value2 = 0.9983391373482872
$ws.Cells.Item(10, 3).Value2 = 0.9732947336262431
$ws.Cells.Item(10, 4).Value2 = 0.99983680146318
$ws.Cells.Item(10, 5).Value2 = 0.9999686768660637
